$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25, pushing the "CONSUMER" feature block
# (old rows 25-30) down to rows 26-31.
$ws.Rows.Item(25).Insert()

# The new row 25 is a sub-feature of the "DEALERS" block (A16:A24) and
# introduces the new feature "PAYMENT".
$ws.Range("B25").Value = "PAYMENT"

# Extend the "DEALERS" merged header cell (A16:A24) down to include the
# newly inserted row (A16:A25).
$ws.Range("A16:A25").Merge()

# Restore the final selection/view state from the saved workbook.
$ws.Range("C21").Select() | Out-Null
